$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The filled-in clearance is being turned back into a blank template:
# a couple of the pre-filled values are replaced with underline blanks
# and the bookmarks that only existed to mark those now-blanked values
# are removed, while the bookmarks that still wrap a (now blank) value
# are kept intact (name / ctcNumber / dateIssued / address /
# dateIssued2).
# ------------------------------------------------------------------

# 1) "Clearance No.: 0038-2021" -> "Clearance No.: _______"
#    Drop the "clearanceNumber" bookmark that wrapped the number - it
#    is no longer referenced anywhere.
$d.Bookmarks.Item("clearanceNumber").Delete()
$d.Content.Find.Execute("0038-2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "_______", 2) | Out-Null

# 2) Community Tax Certificate number "15809537" -> "________".
#    The "ctcNumber" bookmark keeps wrapping the (now blank) value.
$d.Content.Find.Execute("15809537", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "________", 2) | Out-Null

# 3) Tighten the "dateIssued" bookmark so it stops right after the
#    year instead of swallowing the trailing space.
$bmDate = $d.Bookmarks.Item("dateIssued")
$dateStart = $bmDate.Range.Start
$dateEnd = $bmDate.Range.End
$bmDate.Delete()
$d.Bookmarks.Add("dateIssued", $d.Range($dateStart, $dateEnd - 1)) | Out-Null

# 4) "This certification is issued upon his/her request for: KCCDFI
#    LOAN REQUIREMENT." -> blank underline. Drop the "request"
#    bookmark - it is no longer referenced anywhere.
$d.Bookmarks.Item("request").Delete()
$d.Content.Find.Execute("KCCDFI LOAN REQUIREMENT", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "_____________", 2) | Out-Null
